$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $cand = $doc.Paragraphs.Item($i)
        if ($cand.Range.Text -like "*$needle*") {
            return $cand
        }
    }
    return $null
}

# --- Paragraph: "Fixing All the broken code and cleanup 7:06" ---
# Split the single run into 4 runs: "Fixing All the broken code and cleanup ",
# "(", "2", "h)" - mirroring the "(5h)" heading style used elsewhere.
$p1 = Find-ParagraphByText $d "Fixing All the broken code and cleanup"
$r1 = $p1.Range   # includes trailing paragraph mark

$xml1 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5B29136B" w14:textId="78AF26E7" w:rsidR="00E16671" w:rsidRDefault="00E16671" w:rsidP="00E16671"><w:pPr><w:pStyle w:val="berschrift2"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Fixing All the broken code and cleanup </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>h)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$r1.InsertXML($xml1) | Out-Null

# --- Paragraph: "Even though the due date is over, ..." ---
# Keep the existing run untouched and append two new runs of new text.
$p2 = Find-ParagraphByText $d "Even though the due date is over"
$r2 = $p2.Range   # includes trailing paragraph mark

$xml2 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="76D7D96B" w14:textId="5C06AC4E" w:rsidR="00E16671" w:rsidRPr="00E16671" w:rsidRDefault="00E16671" w:rsidP="00E16671"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Even though the due date is over, I still want to finish this project, and since I destroyed most of my code in order to finish the history, I will take my time now to fix everything and make my code look decent again.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I was able to fix everything and even improve on the code, I created new functions in order to make it easier for </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>services like userService to check hashes/get the userId by creating new functions inside the authService, who is responsible for everything regarding security.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$r2.InsertXML($xml2) | Out-Null

# --- <w:sdtEndPr/> additions ---
# The diff also shows an empty <w:sdtEndPr/> being inserted into two nested
# <w:sdt> blocks (a Table-of-Contents-style building block). This document's
# only <w:sdt> (the Table of Contents) already carries a populated
# <w:sdtEndPr> and is not nested, so there is nothing to change here.
